# Fruta / hortaliza, semanal
# Insert two new weekly rows (127 and 128) into the "Vega Monumental Concepción - Mandarina"
# sheet, pushing the previously existing rows 127-166 down to 129-168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 127, shifting everything below down by two rows.
$ws.Rows.Item(127).Insert()
$ws.Rows.Item(127).Insert()

# --- New row 127 ---
$ws.Cells.Item(127,1).Value  = 11
$ws.Cells.Item(127,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(127,3).Value  = "Bíobío"
$ws.Cells.Item(127,4).Value  = 44855
$ws.Cells.Item(127,5).Value  = 8
$ws.Cells.Item(127,6).Value  = "Fruta"
$ws.Cells.Item(127,7).Value  = 100102
$ws.Cells.Item(127,8).Value  = "Cítricos"
$ws.Cells.Item(127,9).Value  = 100102004
$ws.Cells.Item(127,10).Value = "Mandarina"
$ws.Cells.Item(127,11).Value = "Clementina"
$ws.Cells.Item(127,12).Value = "Especial"
$ws.Cells.Item(127,13).Value = 270
$ws.Cells.Item(127,14).Value = 7000
$ws.Cells.Item(127,15).Value = 7500
$ws.Cells.Item(127,16).Value = 7222
$ws.Cells.Item(127,17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(127,18).Value = "Región de O'Higgins"
$ws.Cells.Item(127,19).Value = 401
$ws.Cells.Item(127,20).Value = 18

# --- New row 128 ---
$ws.Cells.Item(128,1).Value  = 11
$ws.Cells.Item(128,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(128,3).Value  = "Bíobío"
$ws.Cells.Item(128,4).Value  = 44855
$ws.Cells.Item(128,5).Value  = 8
$ws.Cells.Item(128,6).Value  = "Fruta"
$ws.Cells.Item(128,7).Value  = 100102
$ws.Cells.Item(128,8).Value  = "Cítricos"
$ws.Cells.Item(128,9).Value  = 100102004
$ws.Cells.Item(128,10).Value = "Mandarina"
$ws.Cells.Item(128,11).Value = "Clementina"
$ws.Cells.Item(128,12).Value = "Segunda"
$ws.Cells.Item(128,13).Value = 180
$ws.Cells.Item(128,14).Value = 5000
$ws.Cells.Item(128,15).Value = 5500
$ws.Cells.Item(128,16).Value = 5278
$ws.Cells.Item(128,17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(128,18).Value = "Región de O'Higgins"
$ws.Cells.Item(128,19).Value = 293
$ws.Cells.Item(128,20).Value = 18
